$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: opus_big Simple aWCE  -> add "Loaded" marker at H15, move selection
# ---------------------------------------------------------------------------
$wsSimple = $wb.Worksheets.Item("opus_big Simple aWCE")
$wsSimple.Activate()
$wsSimple.Range("H15").Value = "Loaded"
$wsSimple.Range("H21").Select()

# ---------------------------------------------------------------------------
# Sheet: opus_big AoN aWCE  -> add "Loaded" marker at H10, move selection
# ---------------------------------------------------------------------------
$wsAoN = $wb.Worksheets.Item("opus_big AoN aWCE")
$wsAoN.Activate()
$wsAoN.Range("H10").Value = "Loaded"
$wsAoN.Range("I10").Select()

# ---------------------------------------------------------------------------
# Sheet: opus_big Fine aWCE -> add "Loaded" marker at H11, move selection
# ---------------------------------------------------------------------------
$wsFine = $wb.Worksheets.Item("opus_big Fine aWCE")
$wsFine.Activate()
$wsFine.Range("H11").Value = "Loaded"
$wsFine.Range("I11").Select()

# ---------------------------------------------------------------------------
# Sheet: opus_big LSP AoN aWCE  -> fill in newly-measured runtime results
# ---------------------------------------------------------------------------
$wsLspAoN = $wb.Worksheets.Item("opus_big LSP AoN aWCE ")
$wsLspAoN.Activate()

# Row 2's Compute column (F2) was left with stale formatting (no top
# border) from before the row above it existed; match the border used by
# the rest of the column (as seen on F3) before writing its value.
$wsLspAoN.Range("F3").Copy()
$wsLspAoN.Range("F2").PasteSpecial(-4122)

$wsLspAoN.Range("C2").Value = 43.351300000000002
$wsLspAoN.Range("D2").Value = 28698.839
$wsLspAoN.Range("F2").Value = 258.44920000000002

$wsLspAoN.Range("C3").Value = 43.152999999999999
$wsLspAoN.Range("D3").Value = 22930.547200000001
$wsLspAoN.Range("F3").Value = 258.44920000000002

$wsLspAoN.Range("C4").Value = 43.043700000000001
$wsLspAoN.Range("D4").Value = 17618.741399999999
$wsLspAoN.Range("F4").Value = 103.4242

$wsLspAoN.Range("C5").Value = 42.890999999999998
$wsLspAoN.Range("D5").Value = 16430.182199999999
$wsLspAoN.Range("F5").Value = 103.4242

$wsLspAoN.Range("C6").Value = 43.021299999999997
$wsLspAoN.Range("D6").Value = 13641.8586
$wsLspAoN.Range("F6").Value = 155.0676

$wsLspAoN.Range("C10").Value = 43.244900000000001
$wsLspAoN.Range("D10").Value = 17283.948400000001
$wsLspAoN.Range("F10").Value = 103.4242

$wsLspAoN.Range("C11").Value = 43.005299999999998
$wsLspAoN.Range("D11").Value = 16432.360100000002
$wsLspAoN.Range("F11").Value = 155.0676

$wsLspAoN.Range("C13").Value = 43.196599999999997
$wsLspAoN.Range("D13").Value = 28871.741600000001
$wsLspAoN.Range("F13").Value = 258.44920000000002

$wsLspAoN.Range("D24").Select()

# ---------------------------------------------------------------------------
# Sheet: opus_big LSP Fine aWCE -> fill in newly-measured runtime results
# ---------------------------------------------------------------------------
$wsLspFine = $wb.Worksheets.Item("opus_big LSP Fine aWCE ")
$wsLspFine.Activate()

$wsLspFine.Range("C10").Value = 42.866300000000003
$wsLspFine.Range("D10").Value = 17490.058499999999
$wsLspFine.Range("F10").Value = 103.4242

$wsLspFine.Range("C14").Value = 42.865099999999998
$wsLspFine.Range("D14").Value = 15471.5005
$wsLspFine.Range("F14").Value = 103.4242

$wsLspFine.Range("E23").Select()
